# Slide 11 ("ICL results complete for similarity"): the caption text box
# ("TextBox 1") currently holds its text split across two runs:
#   Run 1: "Confident in GPT2 vs Mistral ICL "
#   Run 2: "Natural Instructions"
# Merge them into a single run (keeping the formatting/rPr of the first
# run) so the text reads as one continuous run:
#   "Confident in GPT2 vs Mistral ICL Natural Instructions"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)
$shp = $s.Shapes.Item("TextBox 1")

$finalText = "Confident in GPT2 vs Mistral ICL Natural Instructions"

# The getter already reports the two runs concatenated, so assigning the
# same final string directly would look like a no-op to the host and
# leave the two separate <a:r> runs untouched. Stage a genuinely
# different value first so the subsequent assignment actually rewrites
# the paragraph into a single run.
$shp.TextFrame.TextRange.Text = "placeholder"
$shp.TextFrame.TextRange.Text = $finalText
